$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised historical values in columns C-F (open/high/low/close all equal)
$ws.Range("C358:F358").Value = 30734900000000
$ws.Range("C359:F359").Value = 31322400000000
$ws.Range("C361:F361").Value = 33495800000000
$ws.Range("C362:F362").Value = 32988400000000
$ws.Range("C363:F363").Value = 33552400000000
$ws.Range("C364:F364").Value = 33970200000000
$ws.Range("C365:F365").Value = 33881500000000
$ws.Range("C366:F366").Value = 33795100000000
$ws.Range("C367:F367").Value = 34176400000000
$ws.Range("C368:F368").Value = 34635200000000
$ws.Range("C369:F369").Value = 34976400000000
$ws.Range("C370:F370").Value = 35509700000000
$ws.Range("C371:F371").Value = 36350700000000
$ws.Range("C372:F372").Value = 37653400000000
$ws.Range("C373:F373").Value = 38869800000000
$ws.Range("C375:F375").Value = 39527800000000
$ws.Range("C376:F376").Value = 39600300000000
$ws.Range("C377:F377").Value = 39874100000000
$ws.Range("C380:F380").Value = 41209500000000
$ws.Range("C381:F381").Value = 42040500000000
$ws.Range("C382:F382").Value = 41918300000000
$ws.Range("C387:F387").Value = 40513200000000

# Append two new monthly rows (388-389), copying row 387's date-column
# formatting (style index 2: centered/top, bordered, YYYY-MM-DD HH:MM:SS)
# onto the new A388:A389 cells before filling in values.
$ws.Range("A387").Copy()
$ws.Range("A388:A389").PasteSpecial(-4122)

$ws.Range("A388").Value = 44986.45833333334
$ws.Range("B388").Value = "ECONOMICS:HUM2"
$ws.Range("C388:F388").Value = 40410200000000
$ws.Range("G388").Value = 0

$ws.Range("A389").Value = 45017.45833333334
$ws.Range("B389").Value = "ECONOMICS:HUM2"
$ws.Range("C389:F389").Value = 39781400000000
$ws.Range("G389").Value = 0
